# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# - Updates the "VALOR MORA" total (E11) and "Cant. Periodos" count (F13)
# - Adds a new EC (estado de cuenta) period row "2509" for the existing
#   worker (CC 45530426 - ANA DEL CARMEN CONTRERAS HERRERA), inserted right
#   after the last existing data row, reusing that row's values/format and
#   moving the bottom-border row styling down to the new last data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the account summary values.
$ws.Range("E11").Value = 286652
$ws.Range("F13").Value = 5

# Insert a new row below the last detail row (row 19) so rows below
# (the blank spacer rows + the signature block) shift down by one.
$ws.Rows("20:20").Insert()

# Duplicate row 19 (values + formatting) into the freshly inserted row 20 -
# this is the new "2509" period line for the same worker.
$ws.Range("B19:J19").Copy($ws.Range("B20:J20"))

# Row 19 is no longer the last row in the table, so it should pick up the
# "interior" row formatting (same as rows 16-18) instead of the bottom
# border it had before; row 20 keeps the bottom-border styling that row 19
# used to have.
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# Set the new row's period value.
$ws.Range("E20").Value = "2509"
